$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = "SU-T76"
$ws.Range("H5").Value = "village_details.xlsx"

# Row 6
$ws.Range("A6").Value = "SU-T77"
$ws.Range("B6").Value = "Ujjain"
$ws.Range("C6").Value = "Maksi"
$ws.Range("D6").Value = "Ktahit"
$ws.Range("E6").Value = 40
$ws.Range("F6").Value = "Ankit"
$ws.Range("G6").Value = "Contact Number"
$ws.Range("H6").Value = "village_details (1).xlsx.crdownload"

# Row 7
$ws.Range("A7").Value = "SU-T81"
$ws.Range("E7").Value = 55

# Row 8
$ws.Range("A8").Value = "SU-T1084"
$ws.Range("C8").Value = "Sehore"

# Column width adjustments (matching resulting autofit-like widths in the target file)
$ws.Columns.Item(8).ColumnWidth = 29.6328125
$ws.Columns.Item(9).ColumnWidth = 14.54296875

# Update selection to match final state
$ws.Range("B8").Select() | Out-Null
